$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-01-29 Thursday" "2026-01-30 Friday"

Replace-Text "633×4=2532" "205×5=1025"
Replace-Text "325×9=2925" "196×5=980"
Replace-Text "176×7=1232" "239×9=2151"
Replace-Text "345×8=2760" "320×8=2560"
Replace-Text "176×6=1056" "265×6=1590"

Replace-Text "854×8=6832" "912×8=7296"
Replace-Text "234×5=1170" "868×8=6944"
Replace-Text "452×9=4068" "202×5=1010"
Replace-Text "423×7=2961" "377×3=1131"
Replace-Text "669×4=2676" "430×6=2580"

Replace-Text "727×6=4362" "402×3=1206"
Replace-Text "188×4=752" "908×9=8172"
Replace-Text "668×6=4008" "512×6=3072"
Replace-Text "735×3=2205" "310×9=2790"
Replace-Text "554×2=1108" "764×5=3820"

Replace-Text "713×2=1426" "906×7=6342"
Replace-Text "139×2=278" "201×7=1407"
Replace-Text "190×9=1710" "360×9=3240"
Replace-Text "504×4=2016" "502×4=2008"
Replace-Text "866×6=5196" "396×2=792"

Replace-Text "635×9=5715" "428×4=1712"
Replace-Text "873×9=7857" "862×3=2586"
Replace-Text "552×7=3864" "661×9=5949"
Replace-Text "499×6=2994" "565×6=3390"
Replace-Text "586×7=4102" "372×8=2976"

Write-Output "Done applying replacements"
